$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Treatment query (row 5 / "TreatmentTab"): remove the redundant
# CONCAT() wrapper around REPLACE() in the "Treatment Agent" column.
$old = $ws.Range("B5").Value2
$new = $old.Replace(
    "CONCAT(REPLACE(trt.treatment_agent, ';', ', ')) AS ""Treatment Agent""",
    "REPLACE(trt.treatment_agent, ';', ', ') AS ""Treatment Agent"""
)
$ws.Range("B5").Value2 = $new

# Re-apply the same visual formatting (12pt, theme text color, wrap text)
# that the query cells already use, so the cell keeps its normal look.
$ws.Range("B5").WrapText = $true
$ws.Range("B5").Font.Size = 12
$ws.Range("B5").Font.ThemeColor = 1

# Update the saved selection/view to rest on C5 (as in the edited workbook).
$ws.Range("C5").Select()
